$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as text, preserving its original style,
# to prevent Excel auto-converting numeric-looking strings (e.g. "1.00", "0.489")
# into actual numbers.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '68.385.19'
$ws.Range('D3').Value = '3.599.95'
$ws.Range('E3').Value = '  -2.48%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue $ws.Range('D5') '623.22'
$ws.Range('E5').Value = '  -7.15%  '
Set-TextValue $ws.Range('D6') '156.09'
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('D7').Value = '3.597.54'
$ws.Range('E7').Value = '  -2.46%  '
$ws.Range('E8').Value = '  +0.04%  '
Set-TextValue $ws.Range('D9') '0.489'
$ws.Range('E9').Value = '  -2.32%  '
$ws.Range('E10').Value = '  -2.98%  '
Set-TextValue $ws.Range('D11') '6.95'
$ws.Range('E11').Value = '  -2.29%  '
Set-TextValue $ws.Range('D12') '0.434'
$ws.Range('E12').Value = '  -1.83%  '
Set-TextValue $ws.Range('D13') '0.0000225'
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').Value = '4.212.34'
$ws.Range('E14').Value = '  -2.37%  '
Set-TextValue $ws.Range('D15') '32.13'
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').Value = '3.622.27'
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = '68.353.06'
$ws.Range('E17').Value = '  -1.93%  '
Set-TextValue $ws.Range('D18') '0.117'
$ws.Range('E18').Value = '  +0.44%  '
Set-TextValue $ws.Range('D19') '6.44'
$ws.Range('E19').Value = '  -0.58%  '
Set-TextValue $ws.Range('D20') '15.65'
$ws.Range('E20').Value = '  -3.24%  '
Set-TextValue $ws.Range('D21') '460.48'
$ws.Range('E21').Value = '  -2.27%  '
Set-TextValue $ws.Range('D22') '9.82'
$ws.Range('E22').Value = '  +0.54%  '
Set-TextValue $ws.Range('D23') '0.642'
$ws.Range('E23').Value = '  -0.92%  '
Set-TextValue $ws.Range('D24') '78.16'
$ws.Range('E24').Value = '  -2.16%  '
$ws.Range('D25').Value = '3.746.36'
$ws.Range('E25').Value = '  -2.40%  '
$ws.Range('E26').Value = '  +0.01%  '
Set-TextValue $ws.Range('D27') '10.75'
$ws.Range('E27').Value = '  -2.11%  '
$ws.Range('E28').Value = '  -8.48%  '
Set-TextValue $ws.Range('D29') '8.42'
$ws.Range('E29').Value = '  -7.32%  '
Set-TextValue $ws.Range('D30') '2.59'
$ws.Range('E30').Value = '  -3.81%  '
Set-TextValue $ws.Range('D31') '1.65'
$ws.Range('E31').Value = '  -3.69%  '
Set-TextValue $ws.Range('D32') '1.00'
$ws.Range('E32').Value = '  -0.03%  '
Set-TextValue $ws.Range('D33') '26.24'
$ws.Range('E33').Value = '  -2.26%  '
Set-TextValue $ws.Range('D34') '1.92'
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value = '3.601.36'
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D36') '0.159'
$ws.Range('E36').Value = '  -4.99%  '
Set-TextValue $ws.Range('D37') '6.20'
$ws.Range('E37').Value = '  -4.75%  '
Set-TextValue $ws.Range('D38') '8.14'
$ws.Range('E38').Value = '  -4.37%  '
$ws.Range('E39').Value = '  +0.05%  '
Set-TextValue $ws.Range('D40') '177.79'
$ws.Range('E40').Value = '  +0.43%  '
Set-TextValue $ws.Range('D41') '1.00'
$ws.Range('E41').Value = '  +0.01%  '
Set-TextValue $ws.Range('D42') '5.64'
$ws.Range('E42').Value = '  -7.92%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D43') '0.0882'
$ws.Range('E43').Value = '  -3.10%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D44') '2.15'
$ws.Range('E44').Value = '  -5.31%  '
$ws.Range('E45').Value = '  -3.48%  '
Set-TextValue $ws.Range('D46') '46.14'
$ws.Range('E46').Value = '  -1.89%  '
Set-TextValue $ws.Range('D47') '28.38'
$ws.Range('E47').Value = '  +2.95%  '
Set-TextValue $ws.Range('D48') '2.61'
$ws.Range('E48').Value = '  -5.46%  '
Set-TextValue $ws.Range('D49') '7.73'
$ws.Range('E49').Value = '  -2.09%  '
$ws.Range('E50').Value = '  -7.07%  '
Set-TextValue $ws.Range('D51') '1.01'
$ws.Range('E51').Value = '  -6.19%  '
